# feat: add 2022-Q4 data
#
# The existing "2022-Q3" sheet becomes "2022-Q4" (its data is replaced with
# the new quarter's fund table); a fresh "2022-Q3" sheet is inserted right
# after it holding a verbatim copy of what used to be in the original
# "2022-Q3" sheet. The "总计" (totals) sheet gets a new row for 2022-Q4,
# pushing the existing 2022-Q3 total row down one row.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsOldQ3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Duplicate the current "2022-Q3" sheet so its original data survives
#    verbatim in a new tab placed right after it.
# ---------------------------------------------------------------------
$wsOldQ3.Copy($null, $wsOldQ3)
$wsNewQ3 = $wb.Worksheets.Item(3)

# 2) Rename the original sheet (it now holds the 2022-Q4 figures), then
#    rename the duplicate back to "2022-Q3".
$wsQ4 = $wsOldQ3
$wsQ4.Name = "2022-Q4"
$wsNewQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 3) The new table only has 5 funds (6 rows incl. header) vs. 6 funds
#    before, so drop the now-stale 7th row.
# ---------------------------------------------------------------------
$wsQ4.Rows.Item(7).Clear()

# 4) Re-apply the header/index-column style (bordered+bold, style index 2
#    on the totals sheet) without minting a new style entry.
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsQ4.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Helper values are written in two ways:
#  - True numbers (index column, rank column) via plain .Value assignment.
#  - Numeric-looking text (fund code, scale/position percentages) via a
#    text formula ("=""014271""") then Copy/PasteSpecial-values, so Excel
#    keeps them as text instead of auto-converting to numbers (and no
#    NumberFormat/quote-prefix style gets minted along the way).
# ---------------------------------------------------------------------

# Header labels.
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Index column + rank column (real numbers).
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("A5").Value = 3
$wsQ4.Range("A6").Value = 4
$wsQ4.Range("H2").Value = 4
$wsQ4.Range("H3").Value = 4
$wsQ4.Range("H4").Value = 9
$wsQ4.Range("H5").Value = 4
$wsQ4.Range("H6").Value = 4

# Fund names (already non-numeric text, safe to assign directly).
$wsQ4.Range("C2").Value = "大成北交所两年定开混合A"
$wsQ4.Range("C3").Value = "广发北交所精选两年定开混合A"
$wsQ4.Range("C4").Value = "华夏北交所创新中小企业精选两年定开混合"
$wsQ4.Range("C5").Value = "大成北交所两年定开混合C"
$wsQ4.Range("C6").Value = "广发北交所精选两年定开混合C"

# Numeric-looking text columns: fund code (B) + D/E/F/G figures.
$wsQ4.Range("B2").Formula = '="014271"'
$wsQ4.Range("D2").Formula = '="3.24"'
$wsQ4.Range("E2").Formula = '="68.93"'
$wsQ4.Range("F2").Formula = '="7.34"'
$wsQ4.Range("G2").Formula = '="0.2378"'

$wsQ4.Range("B3").Formula = '="014273"'
$wsQ4.Range("D3").Formula = '="3.23"'
$wsQ4.Range("E3").Formula = '="83.79"'
$wsQ4.Range("F3").Formula = '="6.26"'
$wsQ4.Range("G3").Formula = '="0.2022"'

$wsQ4.Range("B4").Formula = '="014283"'
$wsQ4.Range("D4").Formula = '="3.27"'
$wsQ4.Range("E4").Formula = '="90.95"'
$wsQ4.Range("F4").Formula = '="3.91"'
$wsQ4.Range("G4").Formula = '="0.1279"'

$wsQ4.Range("B5").Formula = '="014272"'
$wsQ4.Range("D5").Formula = '="0.77"'
$wsQ4.Range("E5").Formula = '="68.93"'
$wsQ4.Range("F5").Formula = '="7.34"'
$wsQ4.Range("G5").Formula = '="0.0565"'

$wsQ4.Range("B6").Formula = '="014274"'
$wsQ4.Range("D6").Formula = '="0.81"'
$wsQ4.Range("E6").Formula = '="83.79"'
$wsQ4.Range("F6").Formula = '="6.26"'
$wsQ4.Range("G6").Formula = '="0.0507"'

# Freeze all the text formulas above into plain text values in one shot.
$wsQ4.Range("B2:B6").Copy()
$wsQ4.Range("B2:B6").PasteSpecial(-4163)
$wsQ4.Range("D2:G6").Copy()
$wsQ4.Range("D2:G6").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 5) Totals sheet: push the existing 2022-Q3 row down to row 3 (keeping
#    its bordered index-cell style), then write the new 2022-Q4 row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 6
$wsTotal.Range("D3").Value = 0.66

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.68
